# Apply the "all insert queries complete" edit to the World Cup sample
# ticket table: fill in the previously-blank TICKET# (col J) and
# PRICE (col L) columns for rows 11-15, format the new PRICE values as
# currency, fix the "Goalkeeper" role label on row 28, and update the
# sheet view (zoom + selection) to match where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ticket # (column J) --------------------------------------------------
$ws.Range("J11").Value = 27
$ws.Range("J12").Value = 33
$ws.Range("J13").Value = 59
$ws.Range("J14").Value = 11
$ws.Range("J15").Value = 96

# --- Price (column L) ------------------------------------------------------
$ws.Range("L11").Value = 19.99
$ws.Range("L12").Value = 20
$ws.Range("L13").Value = 20.01
$ws.Range("L14").Value = 1000
$ws.Range("L15").Value = 5

# Apply a currency number format to the new price column as a block so the
# workbook gets a single shared style for it.
$ws.Range("L11:L15").NumberFormat = "#,##0.00"

# --- Fix role label on row 28 (GoalKeeper -> Goalkeeper) -------------------
$ws.Range("N28").Value = "Goalkeeper"

# --- View state: zoom to 200% and select the range being worked on ---------
$excel.ActiveWindow.Zoom = 200
$ws.Range("E19:G28").Select() | Out-Null
